# "Enabled hevc for P1" (confirmed with MediaInfo)
#
# TO DO sheet: rows are feature checks, columns B..F are Pixel 1..Pixel 4.
# Pixel 1 (column B) values are updated for several rows:
#   - Row 13 Synthetic Fill Flash : Pixel1 "Y" -> cleared (blank)
#   - Row 14 Motion Photos        : Pixel1 blank -> "N/A"
#   - Row 15 Top Shot             : Pixel1 blank -> "N/A"
#   - Row 16 HEVC                 : Pixel1 blank -> "Y"   (the actual fix)
#   - Row 24 Tracking AF Night Sight : Pixel1 blank -> "Y"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 - Synthetic Fill Flash: Pixel 1 no longer "Y" (cleared)
$ws.Range("B13").Value = ""

# Row 14 - Motion Photos: Pixel 1 is "N/A" (grey, non-italic style, same as B12)
$ws.Range("B14").Value = "N/A"
$ws.Range("B14").Font.Color = 10066329

# Row 15 - Top Shot: Pixel 1 is "N/A"
$ws.Range("B15").Value = "N/A"
$ws.Range("B15").Font.Color = 10066329

# Row 16 - HEVC: Pixel 1 is now "Y" (blue "Y" style, same as other Y cells)
$ws.Range("B16").Value = "Y"
$ws.Range("B16").Font.Color = 12611584

# Row 24 - Tracking AF Night Sight: Pixel 1 is now "Y"
$ws.Range("B24").Value = "Y"

# Matches the cursor/selection left behind in the saved workbook
$ws.Range("H38").Select() | Out-Null
